# Updates cryptos list prices/volumes (and reorders the Fetch.AI / dogwifhat
# rows) to match the GitHub Actions data refresh.
# Note: numeric-looking price strings are written with a leading "'" so Excel
# keeps them as text (matching the original inlineStr cells) instead of
# auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.118.83'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '3.582.68'
$ws.Range("E3").Value = '  -1.24%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'578.99"
$ws.Range("E5").Value = '  -2.37%  '

$ws.Range("D6").Value = "'188.52"
$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("E7").Value = '  -2.52%  '

$ws.Range("D8").Value = '3.577.81'
$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("E10").Value = '  -1.97%  '

$ws.Range("E11").Value = '  -0.59%  '

$ws.Range("D12").Value = "'56.04"
$ws.Range("E12").Value = '  -3.50%  '

$ws.Range("E13").Value = '  +1.38%  '

$ws.Range("D14").Value = "'9.61"
$ws.Range("E14").Value = '  -1.47%  '

$ws.Range("D15").Value = '4.155.68'
$ws.Range("E15").Value = '  -1.37%  '

$ws.Range("D16").Value = "'19.89"
$ws.Range("E16").Value = '  +2.71%  '

$ws.Range("D17").Value = '3.581.93'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("D18").Value = '69.994.81'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").Value = "'12.57"
$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("E21").Value = '  -1.25%  '

$ws.Range("D22").Value = "'474.00"
$ws.Range("E22").Value = '  -4.55%  '

$ws.Range("D23").Value = "'19.04"
$ws.Range("E23").Value = '  +13.55%  '

$ws.Range("D24").Value = "'5.05"
$ws.Range("E24").Value = '  -8.58%  '

$ws.Range("D25").Value = "'4.35"
$ws.Range("E25").Value = '  -2.23%  '

$ws.Range("D26").Value = "'88.77"
$ws.Range("E26").Value = '  -2.38%  '

$ws.Range("D27").Value = "'3.04"
$ws.Range("E27").Value = '  -2.40%  '

$ws.Range("D28").Value = "'11.06"
$ws.Range("E28").Value = '  -1.58%  '

$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = '  -0.56%  '

$ws.Range("D30").Value = "'32.12"
$ws.Range("E30").Value = '  -0.76%  '

$ws.Range("D31").Value = "'7.69"
$ws.Range("E31").Value = '  +2.06%  '

$ws.Range("E32").Value = '  +3.03%  '

$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("D34").Value = "'65.86"
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("D35").Value = "'580.34"
$ws.Range("E35").Value = '  -6.22%  '

$ws.Range("D36").Value = "'38.91"
$ws.Range("E36").Value = '  +2.20%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '0.0₃0797'
$ws.Range("E38").Value = '  -4.48%  '

$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("E40").Value = '  -5.77%  '

$ws.Range("E41").Value = '  -5.34%  '

$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = "'2.91"
$ws.Range("E42").Value = '  +8.43%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = "'3.21"
$ws.Range("E43").Value = '  +15.47%  '

$ws.Range("D44").Value = '3.232.53'
$ws.Range("E44").Value = '  -3.22%  '

$ws.Range("D45").Value = "'3.12"
$ws.Range("E45").Value = '  +1.26%  '

$ws.Range("D46").Value = "'0.0440"
$ws.Range("E46").Value = '  -1.33%  '

$ws.Range("D47").Value = "'9.55"
$ws.Range("E47").Value = '  +4.73%  '

$ws.Range("E48").Value = '  +0.77%  '

$ws.Range("E49").Value = '  -0.61%  '

$ws.Range("E50").Value = '  +0.00%  '

$ws.Range("E51").Value = '  -4.10%  '
